$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Step 1: rows 5-11 (existing "erledigt" backlog items) gain Status "erledigt" on rows 10/11 ---
$ws.Range("H10").Value = 'erledigt'
$ws.Range("H11").Value = 'erledigt'

# --- Step 2: make room so the "internal tasks" block (old rows 22-28) moves down to rows 41-47 ---
$ws.Rows("22:40").Insert()

# --- Step 3: rewrite the backlog item rows 12-21 with the updated/reordered content ---
$ws.Range("A12").Value = '5-1'
$ws.Range("B12").Value = 'Als <MM> will ich bereits angelegte Filter auch löschen können.'
$ws.Range("C12").Value = 1
$ws.Range("D12").ClearContents()
$ws.Range("E12").Value = 30
$ws.Range("F12").ClearContents()
$ws.Range("G12").ClearContents()
$ws.Range("H12").ClearContents()

$ws.Range("A13").Value = '3-3'
$ws.Range("B13").Value = 'Als <MM> will ich eigene Keywords in beliebiger Sprache für die positiv / negativ Semtimentanalyse bestimmen können.'
$ws.Range("C13").Value = 10
$ws.Range("D13").ClearContents()
$ws.Range("E13").Value = 27
$ws.Range("F13").ClearContents()
$ws.Range("G13").ClearContents()
$ws.Range("H13").ClearContents()

$ws.Range("A14").Value = '2-1'
$ws.Range("B14").Value = 'Als <MM> möchte ich eine Mindestanzahl von zu sammelden Tweets pro Veranstaltung angeben, damit ich bewerten kann, ob das Analyseergebnis sinnvoll ist.'
$ws.Range("C14").Value = 2
$ws.Range("D14").ClearContents()
$ws.Range("E14").Value = 26
$ws.Range("F14").ClearContents()
$ws.Range("G14").ClearContents()
$ws.Range("H14").ClearContents()

$ws.Range("A15").Value = '4-3'
$ws.Range("B15").Value = 'Als <MM> will ich einzelene Tweets manuell löschen, sodass diese nicht wieder analysiert werden.'
$ws.Range("C15").Value = 2
$ws.Range("D15").ClearContents()
$ws.Range("E15").Value = 20
$ws.Range("F15").Value = 'medium'
$ws.Range("G15").Value = 'Funktionalität: Löschen'
$ws.Range("H15").ClearContents()

$ws.Range("A16").Value = '4-4'
$ws.Range("B16").Value = 'Als <MM> will ich das aktuelle Filterobjekt speichern, sodass ich diese später erneut anwenden kann'
$ws.Range("C16").ClearContents()
$ws.Range("D16").ClearContents()
$ws.Range("E16").Value = 15
$ws.Range("F16").Value = 'high'
$ws.Range("G16").Value = 'Filterobjektspeicherung'
$ws.Range("H16").Value = 'erledigt'

$ws.Range("A17").Value = '4-2'
$ws.Range("B17").Value = 'Als <MM> will ich einzelene Tweets manuell für einen Filter ausblenden, sodass diese nicht analysiert werden.'
$ws.Range("C17").Value = 7
$ws.Range("D17").ClearContents()
$ws.Range("E17").Value = 19
$ws.Range("F17").Value = 'medium'
$ws.Range("G17").Value = 'Funktionalität: Ausblenden'
$ws.Range("H17").ClearContents()

$ws.Range("A18").Value = '5'
$ws.Range("B18").Value = 'Als <MM> will ich zu einer bestimmten Veranstaltung gehörende Daten als .csv exportieren.'
$ws.Range("C18").ClearContents()
$ws.Range("D18").ClearContents()
$ws.Range("E18").Value = 15
$ws.Range("F18").Value = 'low'
$ws.Range("G18").Value = 'Export'
$ws.Range("H18").ClearContents()

$ws.Range("A19").Value = '3-2'
$ws.Range("B19").Value = 'Als <MM> will ich verschiedene Darstellungen der analysierten Daten.'
$ws.Range("C19").Value = 5
$ws.Range("D19").ClearContents()
$ws.Range("E19").Value = 14
$ws.Range("F19").Value = 'high'
$ws.Range("G19").Value = 'Darstellung'
$ws.Range("H19").ClearContents()

$ws.Range("A20").Value = '1-5'
$ws.Range("B20").Value = 'Als <MM> will ich über Tweets informiert werden.'
$ws.Range("C20").ClearContents()
$ws.Range("D20").ClearContents()
$ws.Range("E20").Value = 10
$ws.Range("F20").Value = 'low'
$ws.Range("G20").Value = 'Benachrichtigung'
$ws.Range("H20").ClearContents()

$ws.Range("A21").Value = '1-5'
$ws.Range("B21").Value = 'Als <MM> will ich vorhandene Veranstaltungen klonen, um bestimme Parameter zu übernehmen, aber die ursprüngliche Veranstaltung nicht zu verändern.'
$ws.Range("C21").ClearContents()
$ws.Range("D21").ClearContents()
$ws.Range("E21").Value = 1
$ws.Range("F21").ClearContents()
$ws.Range("G21").ClearContents()
$ws.Range("H21").ClearContents()

# --- Step 4: rows 22/23 are now just blank spacer rows ---
$ws.Range("A22").Value = ''
$ws.Range("B22").Value = ''
$ws.Range("A23").Value = ''
$ws.Range("B23").Value = ''

# --- Step 5: re-apply the AutoFilter over the new data extent, filtering rows where Status (H) is blank ---
$ws.Range("A4:H21").AutoFilter(8, @(""), 7)

# --- Step 6: keep the _FilterDatabase defined name's range in sync with the new filter extent ---
foreach ($n in $wb.Names) {
    if ($n.Name -like "*_FilterDatabase*") {
        $n.RefersTo = "='Product Backlog'!`$A`$4:`$H`$21"
    }
}

# --- Step 7: selection / view bookkeeping to mirror the saved workbook state ---
$ws.Range("B23:B24").Select()
